$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 12599.4
$ws.Range("I69").Value = 10999.5
$ws.Range("J69").Value = 12999.375
$ws.Range("K69").Value = 32998.5
$ws.Range("L69").Value = 38998.125
$ws.Range("M69").Value = -32124.5
$ws.Range("N69").Value = -40746.125
$ws.Range("H72").Value = 12599.4
$ws.Range("I72").Value = 10999.5
$ws.Range("J72").Value = 12999.375
$ws.Range("K72").Value = 98995.5
$ws.Range("L72").Value = 116994.375
$ws.Range("M72").Value = -94627.5
$ws.Range("N72").Value = -125730.375
$ws.Range("H80").Value = 30600.4
$ws.Range("I80").Value = 75295
$ws.Range("J80").Value = 804
$ws.Range("K80").Value = 225885
$ws.Range("L80").Value = 2412
$ws.Range("M80").Value = -224887
$ws.Range("N80").Value = -4408
$ws.Range("H83").Value = 30600.4
$ws.Range("I83").Value = 75295
$ws.Range("J83").Value = 804
$ws.Range("K83").Value = 677655
$ws.Range("L83").Value = 7236
$ws.Range("M83").Value = -672663
$ws.Range("N83").Value = -17220
$ws.Range("H100").Value = 3082.8635
$ws.Range("I100").Value = 2220
$ws.Range("K100").Value = 2220
$ws.Range("M100").Value = -1679
$ws.Range("H106").Value = 5516.769
$ws.Range("I106").Value = 3476.5
$ws.Range("K106").Value = 3476.5
$ws.Range("M106").Value = -2845.5
$ws.Range("H132").Value = 3905.6072
$ws.Range("I132").Value = 3037.1924
$ws.Range("J132").Value = 15195
$ws.Range("K132").Value = 9111.5772
$ws.Range("L132").Value = 45585
$ws.Range("M132").Value = -6581.5772
$ws.Range("N132").Value = -50645
$ws.Range("H137").Value = 1528.579
$ws.Range("I137").Value = 912.53845
$ws.Range("J137").Value = 2863.3333
$ws.Range("K137").Value = 2737.61535
$ws.Range("L137").Value = 8589.999899999999
$ws.Range("M137").Value = -187.61535
$ws.Range("N137").Value = -13689.9999
$ws.Range("H140").Value = 83602.5
$ws.Range("J140").Value = 96759.89999999999
$ws.Range("L140").Value = 96759.89999999999
$ws.Range("N140").Value = -107119.9

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 35000
$ws.Range("J24").Value = 35000
$ws.Range("L24").Value = 35000
$ws.Range("N24").Value = -35748
$ws.Range("H32").Value = 23611.04
$ws.Range("I32").Value = 26074.809
$ws.Range("J32").Value = 10676.25
$ws.Range("K32").Value = 26074.809
$ws.Range("L32").Value = 10676.25
$ws.Range("M32").Value = -25787.809
$ws.Range("N32").Value = -11250.25
$ws.Range("H88").Value = 2627.2222
$ws.Range("H91").Value = 2627.2222
$ws.Range("H100").Value = 35000
$ws.Range("J100").Value = 35000
$ws.Range("L100").Value = 35000
$ws.Range("N100").Value = -37164

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 69949
$ws.Range("J60").Value = 69966
$ws.Range("L60").Value = 69966
$ws.Range("N60").Value = -71164
$ws.Range("H86").Value = 2208.4375
$ws.Range("I86").Value = 2293.9
$ws.Range("J86").Value = 2066
$ws.Range("K86").Value = 2293.9
$ws.Range("L86").Value = 2066
$ws.Range("M86").Value = -1170.9
$ws.Range("N86").Value = -4312
$ws.Range("H89").Value = 2208.4375
$ws.Range("I89").Value = 2293.9
$ws.Range("J89").Value = 2066
$ws.Range("K89").Value = 11469.5
$ws.Range("L89").Value = 10330
$ws.Range("M89").Value = -5853.5
$ws.Range("N89").Value = -21562
$ws.Range("H100").Value = 27002.715
$ws.Range("J100").Value = 27002.715
$ws.Range("L100").Value = 27002.715
$ws.Range("N100").Value = -29166.715
$ws.Range("H105").Value = 2809.2307
$ws.Range("I105").Value = 2833.7273
$ws.Range("K105").Value = 2833.7273
$ws.Range("M105").Value = -1086.7273

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 53588.125
$ws.Range("I86").Value = 70325.836
$ws.Range("J86").Value = 3375
$ws.Range("K86").Value = 70325.836
$ws.Range("L86").Value = 3375
$ws.Range("M86").Value = -69202.836
$ws.Range("N86").Value = -5621
$ws.Range("H89").Value = 53588.125
$ws.Range("I89").Value = 70325.836
$ws.Range("J89").Value = 3375
$ws.Range("K89").Value = 351629.18
$ws.Range("L89").Value = 16875
$ws.Range("M89").Value = -346013.18
$ws.Range("N89").Value = -28107
$ws.Range("H99").Value = 2605.7693
$ws.Range("I99").Value = 1759.8
$ws.Range("J99").Value = 3134.5
$ws.Range("K99").Value = 1759.8
$ws.Range("L99").Value = 3134.5
$ws.Range("M99").Value = -261.8
$ws.Range("N99").Value = -6130.5
$ws.Range("H126").Value = 2605.7693
$ws.Range("I126").Value = 1759.8
$ws.Range("J126").Value = 3134.5
$ws.Range("K126").Value = 5279.4
$ws.Range("L126").Value = 9403.5
$ws.Range("M126").Value = -2809.4
$ws.Range("N126").Value = -14343.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 13583.167
$ws.Range("I3").Value = 7166.3335
$ws.Range("K3").Value = 21499.0005
$ws.Range("M3").Value = -21387.0005
$ws.Range("H50").Value = 250.9
$ws.Range("I50").Value = 142.5
$ws.Range("K50").Value = 427.5
$ws.Range("M50").Value = 53.5
$ws.Range("H53").Value = 250.9
$ws.Range("I53").Value = 142.5
$ws.Range("K53").Value = 427.5
$ws.Range("M53").Value = 53.5
$ws.Range("H81").Value = 4998.5
$ws.Range("J81").Value = 4998
$ws.Range("L81").Value = 14994
$ws.Range("N81").Value = -17240
$ws.Range("H84").Value = 4998.5
$ws.Range("J84").Value = 4998
$ws.Range("L84").Value = 44982
$ws.Range("N84").Value = -56214
$ws.Range("H131").Value = 2786751.2
$ws.Range("J131").Value = 3232593.2
$ws.Range("L131").Value = 9697779.600000001
$ws.Range("N131").Value = -9707859.600000001
$ws.Range("H133").Value = 10500
$ws.Range("I133").Value = 1000
$ws.Range("K133").Value = 3000
$ws.Range("M133").Value = 2060
$ws.Range("H134").Value = 5866.25
$ws.Range("I134").Value = 5866.25
$ws.Range("K134").Value = 17598.75
$ws.Range("M134").Value = -12528.75

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2150.4075
$ws.Range("I113").Value = 1777.6522
$ws.Range("K113").Value = 1777.6522
$ws.Range("M113").Value = 392.3478
$ws.Range("H126").Value = 7241.1333
$ws.Range("I126").Value = 6653
$ws.Range("J126").Value = 7633.222
$ws.Range("K126").Value = 19959
$ws.Range("L126").Value = 22899.666
$ws.Range("M126").Value = -17489
$ws.Range("N126").Value = -27839.666
$ws.Range("H132").Value = 31619.03
$ws.Range("I132").Value = 38956.15
$ws.Range("K132").Value = 116868.45
$ws.Range("M132").Value = -114338.45

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4516.5713
$ws.Range("I7").Value = 3782.5833
$ws.Range("J7").Value = 5495.222
$ws.Range("K7").Value = 3782.5833
$ws.Range("L7").Value = 5495.222
$ws.Range("M7").Value = -3670.5833
$ws.Range("N7").Value = -5719.222
$ws.Range("H61").Value = 3923.85
$ws.Range("I61").Value = 3693.2222
$ws.Range("K61").Value = 3693.2222
$ws.Range("M61").Value = -3491.2222
$ws.Range("H113").Value = 3923.85
$ws.Range("I113").Value = 3693.2222
$ws.Range("K113").Value = 3693.2222
$ws.Range("M113").Value = -1523.2222
$ws.Range("H126").Value = 4516.5713
$ws.Range("I126").Value = 3782.5833
$ws.Range("J126").Value = 5495.222
$ws.Range("K126").Value = 11347.7499
$ws.Range("L126").Value = 16485.666
$ws.Range("M126").Value = -8877.749899999999
$ws.Range("N126").Value = -21425.666
$ws.Range("H132").Value = 42029.227
$ws.Range("I132").Value = 50812.24
$ws.Range("J132").Value = 5433.3335
$ws.Range("K132").Value = 152436.72
$ws.Range("L132").Value = 16300.0005
$ws.Range("M132").Value = -149906.72
$ws.Range("N132").Value = -21360.0005
$ws.Range("H136").Value = 5443.4
$ws.Range("I136").Value = 5443.4
$ws.Range("K136").Value = 16330.2
$ws.Range("M136").Value = -13780.2

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 14027.2
$ws.Range("I39").Value = 12519.5
$ws.Range("J39").Value = 15032.333
$ws.Range("K39").Value = 12519.5
$ws.Range("L39").Value = 15032.333
$ws.Range("M39").Value = -12106.5
$ws.Range("N39").Value = -15858.333
$ws.Range("H81").Value = 6740.636
$ws.Range("I81").Value = 3477.8572
$ws.Range("J81").Value = 12450.5
$ws.Range("K81").Value = 6955.7144
$ws.Range("L81").Value = 24901
$ws.Range("M81").Value = -5894.7144
$ws.Range("N81").Value = -27023
$ws.Range("H84").Value = 6740.636
$ws.Range("I84").Value = 3477.8572
$ws.Range("J84").Value = 12450.5
$ws.Range("K84").Value = 34778.572
$ws.Range("L84").Value = 124505
$ws.Range("M84").Value = -29474.572
$ws.Range("N84").Value = -135113
